$d = $word.ActiveDocument
$wns = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# ---------------------------------------------------------------------------
# 1) Split the run " equivocada." into " " / "equivocada" (spell-checked) /
#    "." -- matches the proofErr wrapping used elsewhere in the paragraph.
# ---------------------------------------------------------------------------
$find1 = $d.Content
$find1.Find.Execute(" equivocada.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng1 = $d.Range($find1.Start, $find1.End)
$xml1 = '<w:p xmlns:w="' + $wns + '"><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>equivocada</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r></w:p>'
$rng1.InsertXML($xml1)

# ---------------------------------------------------------------------------
# 2) Drop the old _GoBack bookmark (it used to sit right after "equivocada.")
# ---------------------------------------------------------------------------
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# ---------------------------------------------------------------------------
# 3) Insert the new "Instant Gratification" paragraph right after the
#    "equivocada." paragraph, before "The Original Seven Deadly Sins..."
# ---------------------------------------------------------------------------
$find2 = $d.Content
$find2.Find.Execute("equivocada.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$p2 = $find2.Paragraphs(1)
$p2.Range.InsertParagraphAfter()
$newPara = $p2.Next()
$newPara.Range.Text = "Instant Gratification"

# ---------------------------------------------------------------------------
# 4) Re-insert "The Original Seven Deadly Sins- Not so sinful" as its own
#    paragraph, carrying the _GoBack bookmark at its end.
# ---------------------------------------------------------------------------
$sinsPara = $newPara.Next()
$xml4 = '<w:p xmlns:w="' + $wns + '"><w:r><w:t>The Original Seven Deadly Sins- Not so sinful</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$sinsPara.Range.InsertXML($xml4)

# ---------------------------------------------------------------------------
# 5) Move <w:lastRenderedPageBreak/> from the "Binary" run to the
#    "Dichotomy" run (one paragraph earlier).
# ---------------------------------------------------------------------------
$findD = $d.Content
$findD.Find.Execute("Dichotomy", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$pD = $findD.Paragraphs(1)
$rD = $pD.Range
$rDnoMark = $d.Range($rD.Start, $rD.End - 1)
$xmlD = '<w:p xmlns:w="' + $wns + '"><w:r><w:lastRenderedPageBreak/><w:t>Dichotomy</w:t></w:r></w:p>'
$rDnoMark.InsertXML($xmlD)

$findB = $d.Content
$findB.Find.Execute("Binary", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$pB = $findB.Paragraphs(1)
$rB = $pB.Range
$rBnoMark = $d.Range($rB.Start, $rB.End - 1)
$xmlB = '<w:p xmlns:w="' + $wns + '"><w:r><w:t>Binary</w:t></w:r></w:p>'
$rBnoMark.InsertXML($xmlB)

# ---------------------------------------------------------------------------
# 6) Re-split the "However, if the sins..." paragraph's two runs at a
#    different point, without altering the combined text or the position
#    of the <w:lastRenderedPageBreak/> marker (still starts the 2nd run).
# ---------------------------------------------------------------------------
$findH = $d.Content
$findH.Find.Execute("However, if the sins are simply being excessive", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$pH = $findH.Paragraphs(1)
$rH = $pH.Range
$rHnoMark = $d.Range($rH.Start, $rH.End - 1)
$xmlH = '<w:p xmlns:w="' + $wns + '"><w:r><w:t xml:space="preserve">However, if the sins are simply being excessive. What should be considered a cardinal sin? Before, humans were looking into finding how to be a wholesome or “complete” person, full of virtues. This type of belief was held by many individuals regardless of their religion. Which is why I want to move away from the </w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">traditional idea of the “seven deadly sins”, and want to move to a new version of these sins that play a fundamental role in understanding how to be a better person in our modern society. These writings are simply an opinion formed by somebody who wishes for the best. Thus, if there is any disagreement with my thoughts, I believe that everyone can either correct me or contribute to this collection of essays. However, I will never be able to consult with every individual possible before making this public. Thus I thank everyone who played part of this and any constructive criticism will always be welcomed. </w:t></w:r></w:p>'
$rHnoMark.InsertXML($xmlH)
